$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("UI")
$sheet3 = $wb.Worksheets.Item("STR")

# --- Apply wrap/font style (style index 3) to target cells by copying from an existing styled cell ---
$sheet3.Cells.Item(55, 2).Copy($sheet3.Cells.Item(150, 2))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(152, 3))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(155, 3))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(158, 3))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(154, 3))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(156, 3))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(157, 3))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(161, 3))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(164, 3))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(150, 3))
$sheet3.Cells.Item(55, 2).Copy($sheet3.Cells.Item(151, 2))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(151, 3))
$sheet3.Cells.Item(56, 3).Copy($sheet3.Cells.Item(153, 3))

# --- Set new shared-string values in the exact order they were introduced ---
$sheet3.Cells.Item(150, 1).Value = 'STR_PENETRATION'
$sheet3.Cells.Item(150, 2).Value = 'Penetration'
$sheet3.Cells.Item(151, 1).Value = 'STR_CRIT_DMG_RES'
$sheet3.Cells.Item(152, 1).Value = 'STR_ENGLISH'
$sheet3.Cells.Item(152, 3).Value = 'Tiếng Anh'
$sheet3.Cells.Item(153, 1).Value = 'STR_VIETNAMESE'
$sheet3.Cells.Item(154, 1).Value = 'STR_COIN'
$sheet3.Cells.Item(156, 1).Value = 'STR_RELIC_ESSENCE'
$sheet3.Cells.Item(157, 1).Value = 'STR_ARMOR_PRIMORITE'
$sheet3.Cells.Item(154, 2).Value = 'Coin'
$sheet3.Cells.Item(156, 2).Value = 'Relic Essence'
$sheet3.Cells.Item(157, 2).Value = 'Armor primorite'
$sheet3.Cells.Item(158, 1).Value = 'STR_COIN_DES'
$sheet3.Cells.Item(160, 1).Value = 'STR_RELIC_ESSENCE_DES'
$sheet3.Cells.Item(161, 1).Value = 'STR_ARMOR_PRIMORITE_DES'
$sheet3.Cells.Item(159, 1).Value = 'STR_JADE_DES'
$sheet3.Cells.Item(155, 1).Value = 'STR_JADE'
$sheet3.Cells.Item(155, 2).Value = 'Jade'
$sheet3.Cells.Item(158, 2).Value = 'its''s the unified currency issued by the Tang Empire, eccepted everywhere in Earth importal Reaml.'
$sheet3.Cells.Item(159, 2).Value = 'Fashioned from jade, this currency circulates among immortals, buddhas, and even demons.'
$sheet3.Cells.Item(160, 2).Value = 'Ores that contain abundant spiritual energy. Infused with divine power in tis core, the ore''s surface is permeated with spiritual energy.'
$sheet3.Cells.Item(161, 2).Value = 'Crystals that contain abundant spiritual energy.\n Overflowing with spiritual energy, the crystal is forged into solidity by celestial artisans with divine flame.'
$sheet3.Cells.Item(162, 1).Value = 'STR_COIN_USE'
$sheet3.Cells.Item(163, 1).Value = 'STR_JADE_USE'
$sheet3.Cells.Item(164, 1).Value = 'STR_RELIC_ESSENCE_USE'
$sheet3.Cells.Item(165, 1).Value = 'STR_ARMOR_PRIMORITE_USE'
$sheet3.Cells.Item(165, 2).Value = 'Used for armor enhancement.'
$sheet3.Cells.Item(162, 2).Value = 'The most commonly used currency.'
$sheet3.Cells.Item(164, 2).Value = 'Used for Relic level-up.'
$sheet3.Cells.Item(163, 2).Value = 'A special currency circulation arcoss the Three Realms.'
$sheet1.Cells.Item(64, 1).Value = 'UI_SELECT_ENHANCE'
$sheet1.Cells.Item(64, 2).Value = 'Select Enhancement Materials'
$sheet1.Cells.Item(65, 1).Value = 'UI_RECYCLE'
$sheet1.Cells.Item(66, 1).Value = 'UI_TARGET_LEVEL'
$sheet1.Cells.Item(66, 2).Value = 'Target Level'
$sheet1.Cells.Item(65, 2).Value = 'Recycle'
$sheet1.Cells.Item(67, 1).Value = 'UI_ENHANCE_ARMOR'
$sheet1.Cells.Item(67, 2).Value = 'Enhance Armor'
$sheet3.Cells.Item(155, 3).Value = 'Tiên Ngọc'
$sheet3.Cells.Item(159, 3).Value = 'Khắc ngọc thành tiên, viếng thăm thần tiên, gõ cửa phật môn, để vào yêu đạo.'
$sheet3.Cells.Item(158, 3).Value = 'Đồng tiền thống nhất do nhà Đường phát hành có thể lưu thông ở Tiên Giới.'
$sheet3.Cells.Item(163, 3).Value = 'Loại tiền đặc biệt được lưu thông ở Tam Giới.'
$sheet3.Cells.Item(154, 3).Value = 'Thông Bảo'
$sheet3.Cells.Item(162, 3).Value = 'Đá quặng chúa nhiều linh khí.\n Kim loại thần thông, linh hồn xuyên qua bề mặt đến từng chi tiết.'
$sheet3.Cells.Item(156, 3).Value = 'Nguyên Uẩn Pháp Bảo'
$sheet3.Cells.Item(157, 3).Value = 'Nguyên Uẩn Trang Bị'
$sheet3.Cells.Item(161, 3).Value = 'Thạch anh chứa nhiều linh khí.\n Đá thô tràn trề linh khí, được thợ tiên điêu khắc, đốt cháy bằng lừa linh hồn, đông cứng thành pha lê.'
$sheet3.Cells.Item(165, 3).Value = 'Dùng để cường hóa trang bị.'
$sheet3.Cells.Item(164, 3).Value = 'Dùng để cường hóa pháp bảo.'

# --- Set cells that reuse existing shared strings ---
$sheet3.Cells.Item(150, 3).Value = 'Xuyên giáp'
$sheet3.Cells.Item(151, 2).Value = 'Crit DMG Res'
$sheet3.Cells.Item(151, 3).Value = 'Kháng ST bạo kích'
$sheet3.Cells.Item(152, 2).Value = 'ENGLISH'
$sheet3.Cells.Item(153, 2).Value = 'VIETNAMESE'
$sheet3.Cells.Item(153, 3).Value = 'Tiếng Việt'